# Ready to fit incidence.
# Update the incidence2018_plus sheet so that the age-specific incidence
# rates are derived from prevalence2018 using a /2 factor (instead of the
# old /100 factor) for ages 25 and up, and refresh the view/selection
# state to match where the author was last working.

$wb = $excel.ActiveWorkbook

$prevSheet = $wb.Worksheets.Item("prevalence2018")
$incSheet  = $wb.Worksheets.Item("incidence2018_plus")

# --- 1. Fix up the incidence formulas (rows 27-122, column C) ----------
# They used to divide the prevalence by 100; they should divide by 2.
for ($r = 27; $r -le 122; $r++) {
    $incSheet.Range("C$r").Formula = "=prevalence2018!C$r/2"
}

# --- 2. Restore/update the view state for both sheets -------------------
# prevalence2018: no longer the selected tab; scrolled up a bit; whole
# column C selected with the active cell back at the top.
$prevSheet.Activate()
$prevSheet.Range("C1:C1048576").Select()

# incidence2018_plus: now the selected tab, scrolled near the top, with
# a single cell (F16) selected.
$incSheet.Activate()
$incSheet.Range("F16").Select()

Write-Output "Updated incidence2018_plus formulas (C27:C122) and view state"
